$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (name/description boundary).
# This shifts the existing "description" column (D) to column E.
$ws.Columns("D").Insert()

# Populate the new "type" column header and values (row 1 = header).
$ws.Range("D1").Value = "type"
$ws.Range("D2").Value = "integer"
$ws.Range("D3").Value = "string"
$ws.Range("D4").Value = "string"
$ws.Range("D5").Value = "integer"
$ws.Range("D6").Value = "string"
$ws.Range("D7").Value = "string"
$ws.Range("D8").Value = "integer"
$ws.Range("D9").Value = "integer"

# Grow the Excel table (ListObject) to include the new column.
$tbl = $ws.ListObjects.Item(1)
$newCol = $tbl.ListColumns.Add()

# Re-establish the "description" header text so the table picks it up
# as the column name for the shifted column.
$ws.Range("E1").Value = "description"

# Give the new column a sensible width like Excel's auto-fit would.
$ws.Range("D1:D9").EntireColumn.AutoFit()
